$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Insert 6 new rows before row 72 (inside the "Metrics" table), pushing the
# existing "Treatments" rows (old rows 72-74) down to rows 78-80.
$ws.Rows.Item(72).Resize(6).Insert()

# Grow the table definition so the new rows are part of "Metrics" (A1:F74 -> A1:F80)
$tbl.Resize($ws.Range("A1:F80"))

# New "Aged Care" metrics being added (Molnupiravir + Paxlovid prescriptions)
$newMetrics = @(
  "# Aged Care Molnupiravir Prescriptions",
  "# Aged Care Molnupiravir Prescriptions per 1M",
  "% Aged Care Molnupiravir Prescriptions per Case",
  "# Aged Care Paxlovid Prescriptions",
  "# Aged Care Paxlovid Prescriptions per 1M",
  "% Aged Care Paxlovid Prescriptions per Case"
)

for ($i = 0; $i -lt $newMetrics.Length; $i++) {
  $r = 72 + $i
  $ws.Cells.Item($r, 1).Value = "Aged Care"
  $ws.Cells.Item($r, 2).Value = 60
  $ws.Cells.Item($r, 3).Value = $newMetrics[$i]
  $ws.Cells.Item($r, 4).Value = (($r - 1) * 10)
  $ws.Cells.Item($r, 6).Value = "X"
}

# Keep the "Metric - Sort" helper column sequential for the rows that got
# pushed down to make room for the new metrics above.
$ws.Cells.Item(78, 4).Value = 770
$ws.Cells.Item(79, 4).Value = 780
$ws.Cells.Item(80, 4).Value = 790

# Match the saved selection state left behind after entering the new rows.
[void]$ws.Range("F71:F77").Select()
